$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the currency-symbol column (D) to parenthesized display strings
# and the refund percentage column (G); column H (refund amount) and
# column K (the warning message) are formulas and recalculate automatically.
$ws.Cells.Item(2, 4).Value = '($)'
$ws.Cells.Item(2, 7).Value = 13
$ws.Cells.Item(3, 4).Value = '(€)'
$ws.Cells.Item(3, 7).Value = 12.5
$ws.Cells.Item(4, 4).Value = '($)'
$ws.Cells.Item(4, 7).Value = 5.93
$ws.Cells.Item(5, 4).Value = '(€)'
$ws.Cells.Item(5, 7).Value = 11.8
$ws.Cells.Item(6, 4).Value = '(¥)'
$ws.Cells.Item(6, 7).Value = 9
$ws.Cells.Item(7, 4).Value = '(kn)'
$ws.Cells.Item(7, 7).Value = 15
$ws.Cells.Item(8, 4).Value = '(€)'
$ws.Cells.Item(8, 7).Value = 11.3
$ws.Cells.Item(9, 4).Value = '(Kč)'
$ws.Cells.Item(9, 7).Value = 11.7
$ws.Cells.Item(10, 4).Value = '(kr)'
$ws.Cells.Item(10, 7).Value = 13
$ws.Cells.Item(11, 4).Value = '(€)'
$ws.Cells.Item(11, 7).Value = 11.5
$ws.Cells.Item(12, 4).Value = '(€)'
$ws.Cells.Item(12, 7).Value = 13.5
$ws.Cells.Item(13, 4).Value = '(€)'
$ws.Cells.Item(13, 7).Value = 12
$ws.Cells.Item(14, 4).Value = '(€)'
$ws.Cells.Item(14, 7).Value = 11.4
$ws.Cells.Item(15, 4).Value = '(€)'
$ws.Cells.Item(15, 7).Value = 12.4
$ws.Cells.Item(16, 4).Value = '(Ft)'
$ws.Cells.Item(16, 7).Value = 14.258995000000001
$ws.Cells.Item(17, 4).Value = '(kr)'
$ws.Cells.Item(17, 7).Value = 14.166667
$ws.Cells.Item(18, 4).Value = '(€)'
$ws.Cells.Item(18, 7).Value = 13.9
$ws.Cells.Item(19, 4).Value = '(€)'
$ws.Cells.Item(19, 7).Value = 12.9
$ws.Cells.Item(20, 4).Value = '(¥)'
$ws.Cells.Item(20, 7).Value = 6
$ws.Cells.Item(21, 4).Value = '(€)'
$ws.Cells.Item(21, 7).Value = 12.397
$ws.Cells.Item(22, 4).Value = '(£)'
$ws.Cells.Item(22, 7).Value = 8
$ws.Cells.Item(23, 4).Value = '(€)'
$ws.Cells.Item(23, 7).Value = 12.2
$ws.Cells.Item(24, 4).Value = '(€)'
$ws.Cells.Item(24, 7).Value = 9.8000000000000007
$ws.Cells.Item(25, 4).Value = '(RM)'
$ws.Cells.Item(25, 7).Value = 4.8099999999999996
$ws.Cells.Item(26, 4).Value = '(MAD)'
$ws.Cells.Item(26, 7).Value = 12.75
$ws.Cells.Item(27, 4).Value = '(€)'
$ws.Cells.Item(27, 7).Value = 11.2
$ws.Cells.Item(28, 4).Value = '(kr)'
$ws.Cells.Item(28, 7).Value = 12
$ws.Cells.Item(29, 4).Value = '(zł)'
$ws.Cells.Item(29, 7).Value = 13.6
$ws.Cells.Item(30, 4).Value = '(€)'
$ws.Cells.Item(30, 7).Value = 13.9
$ws.Cells.Item(31, 4).Value = '(₩)'
$ws.Cells.Item(31, 7).Value = 5
$ws.Cells.Item(32, 4).Value = '(₽)'
$ws.Cells.Item(32, 7).Value = 11
$ws.Cells.Item(33, 4).Value = '($)'
$ws.Cells.Item(33, 7).Value = 5.5
$ws.Cells.Item(34, 4).Value = '(€)'
$ws.Cells.Item(34, 7).Value = 11.4
$ws.Cells.Item(35, 4).Value = '(€)'
$ws.Cells.Item(35, 7).Value = 14.5
$ws.Cells.Item(36, 4).Value = '(€)'
$ws.Cells.Item(36, 7).Value = 12.85
$ws.Cells.Item(37, 4).Value = '(kr)'
$ws.Cells.Item(37, 7).Value = 11.8
$ws.Cells.Item(38, 4).Value = '(Fr.)'
$ws.Cells.Item(38, 7).Value = 4.5999999999999996
$ws.Cells.Item(39, 4).Value = '(₺)'
$ws.Cells.Item(39, 7).Value = 4.0750000000000002
$ws.Cells.Item(40, 4).Value = '(£)'
$ws.Cells.Item(40, 7).Value = 12.2
$ws.Cells.Item(41, 4).Value = '($U)'
$ws.Cells.Item(41, 7).Value = 14.426

# Update the view: scroll back to the top and select G12 (was C34, scrolled to A13)
$ws.Range("A1").Select() | Out-Null
$ws.Range("G12").Select() | Out-Null
